$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $Address, $Val)
    $range = $Worksheet.Range($Address)
    $range.NumberFormat = "@"
    $range.Value = $Val
    $range.Style = "Normal"
}

Set-TextValue $ws 'D2' '23.879.67'
Set-TextValue $ws 'E2' '  -1.67%  '
Set-TextValue $ws 'D3' '1.650.11'
Set-TextValue $ws 'E3' '  -1.08%  '
Set-TextValue $ws 'E4' '  -0.14%  '
Set-TextValue $ws 'D5' '310.72'
Set-TextValue $ws 'E5' '  -0.56%  '
Set-TextValue $ws 'E6' '  -0.11%  '
Set-TextValue $ws 'D7' '0.3890'
Set-TextValue $ws 'E7' '  -1.85%  '
Set-TextValue $ws 'D8' '0.3834'
Set-TextValue $ws 'E8' '  -2.57%  '
Set-TextValue $ws 'D9' '51.25'
Set-TextValue $ws 'E9' '  -1.31%  '
Set-TextValue $ws 'D10' '1.341'
Set-TextValue $ws 'E10' '  -3.35%  '
Set-TextValue $ws 'D11' '1.001'
Set-TextValue $ws 'E11' '  -0.14%  '
Set-TextValue $ws 'D12' '0.08442'
Set-TextValue $ws 'E12' '  -1.60%  '
Set-TextValue $ws 'D13' '23.84'
Set-TextValue $ws 'E13' '  -2.38%  '
Set-TextValue $ws 'D14' '7.007'
Set-TextValue $ws 'E14' '  -4.24%  '
Set-TextValue $ws 'D15' '8.007'
Set-TextValue $ws 'E15' '  -0.04%  '
Set-TextValue $ws 'E16' '  -1.10%  '
Set-TextValue $ws 'D17' '1.650.34'
Set-TextValue $ws 'E17' '  -0.93%  '
Set-TextValue $ws 'D18' '94.00'
Set-TextValue $ws 'E18' '  -1.84%  '
Set-TextValue $ws 'D19' '0.06974'
Set-TextValue $ws 'D20' '19.53'
Set-TextValue $ws 'E20' '  -4.62%  '
Set-TextValue $ws 'D21' '6.947'
Set-TextValue $ws 'E21' '  -0.82%  '
Set-TextValue $ws 'D22' '1.000'
Set-TextValue $ws 'E22' '  -0.12%  '
Set-TextValue $ws 'D23' '13.62'
Set-TextValue $ws 'E23' '  -1.95%  '
Set-TextValue $ws 'D24' '23.884.55'
Set-TextValue $ws 'E24' '  -1.70%  '
Set-TextValue $ws 'D25' '2.445'
Set-TextValue $ws 'E25' '  -3.73%  '
Set-TextValue $ws 'D26' '2.920'
Set-TextValue $ws 'E26' '  -6.19%  '
Set-TextValue $ws 'D27' '21.94'
Set-TextValue $ws 'E27' '  -2.49%  '
Set-TextValue $ws 'D28' '153.75'
Set-TextValue $ws 'E28' '  -1.95%  '
Set-TextValue $ws 'D29' '5.380'
Set-TextValue $ws 'E29' '  +0.37%  '
Set-TextValue $ws 'D30' '137.23'
Set-TextValue $ws 'E30' '  -3.55%  '
Set-TextValue $ws 'D31' '7.726'
Set-TextValue $ws 'E31' '  -3.51%  '
Set-TextValue $ws 'E32' '  -1.90%  '
Set-TextValue $ws 'D33' '1.831.61'
Set-TextValue $ws 'E33' '  -0.94%  '
Set-TextValue $ws 'D34' '0.08147'
Set-TextValue $ws 'E34' '  -1.46%  '
Set-TextValue $ws 'D35' '0.9893'
Set-TextValue $ws 'E35' '  -6.63%  '
Set-TextValue $ws 'D36' '0.02909'
Set-TextValue $ws 'E36' '  -5.25%  '
Set-TextValue $ws 'D37' '6.657'
Set-TextValue $ws 'E37' '  -3.11%  '
Set-TextValue $ws 'D38' '0.2674'
Set-TextValue $ws 'E38' '  -3.16%  '
Set-TextValue $ws 'D39' '10.50'
Set-TextValue $ws 'E39' '  -5.72%  '
Set-TextValue $ws 'D40' '0.09105'
Set-TextValue $ws 'E40' '  -1.95%  '
Set-TextValue $ws 'D41' '0.7553'
Set-TextValue $ws 'E41' '  -1.68%  '
Set-TextValue $ws 'D42' '13.46'
Set-TextValue $ws 'E42' '  -2.23%  '
Set-TextValue $ws 'D43' '1.421'
Set-TextValue $ws 'E43' '  -1.36%  '
Set-TextValue $ws 'D44' '16.65'
Set-TextValue $ws 'E44' '  +0.41%  '
Set-TextValue $ws 'D45' '0.6934'
Set-TextValue $ws 'E45' '  -1.88%  '
Set-TextValue $ws 'E46' '  -3.43%  '
Set-TextValue $ws 'E47' '  -0.70%  '
Set-TextValue $ws 'D48' '1.0000'
Set-TextValue $ws 'E48' '  -0.10%  '
Set-TextValue $ws 'D49' '0.08271'
Set-TextValue $ws 'E49' '  -1.59%  '
Set-TextValue $ws 'D50' '133.77'
Set-TextValue $ws 'E50' '  -1.99%  '
Set-TextValue $ws 'D51' '1.221'
Set-TextValue $ws 'E51' '  -3.22%  '

Write-Host "Applied 93 cell updates to D2:E51 range"
